$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$ws = $wb.Worksheets.Item("Metadata")

# Version: 1.0.0 -> 1.1.0
$ws.Range("B3").Value = "1.1.0"

# Date: 2023-06-07T11:52:14+02:00 -> 2023-07-10T23:08:03+02:00
$ws.Range("B8").Value = "2023-07-10T23:08:03+02:00"

# --- Style updates: ensure alignment (vertical=top, wrapText) is applied via the
# "applyAlignment" flag on the cell styles used throughout both sheets (s="1" and s="2")
$ws1 = $wb.Worksheets.Item("Metadata")
$ws2 = $wb.Worksheets.Item("Include from SNOMED CT")

foreach ($sheet in @($ws1, $ws2)) {
    $used = $sheet.UsedRange
    $used.VerticalAlignment = -4160  # xlTop
    $used.WrapText = $true
}
